# 16-Apr-2024: GUI implemented to load the configuration file and specify
# how many test paper to be generated.
#
# The config sheet gains two helper columns (C/D) carrying short notes next
# to the "online test" / "paper test" related settings, and the file-name
# values for the test paper / marksheet settings drop their hard-coded file
# extension (the program now appends the extension itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hard-coded file extensions now that the program appends them.
$ws.Range("B11").Value = "testpaper"
$ws.Range("B12").Value = "marksheet"

# New annotation column (C) marking which settings apply to the online test
# vs. the (new) paper test feature.
$ws.Range("C8").Value = "; online test"
$ws.Range("C9").Value = "; online test"
$ws.Range("C10").Value = "; paper test"
$ws.Range("C11").Value = "; paper test"
$ws.Range("C12").Value = "; paper test"

# New annotation column (D) documenting the automatic file-extension
# behaviour for the paper-test file name settings.
$ws.Range("D11").Value = "the file type (pdf) will be added by the program"
$ws.Range("D12").Value = "the file type (xlsx) will be added by the program"

# Clear the stale cell selection (previously B10) left over from editing,
# resetting it back to the sheet's home cell.
$ws.Range("A1").Select()
